$d = $word.ActiveDocument

# Locate the plain-text e-mail address run and turn it into a mailto: hyperlink,
# matching the direct character formatting already used for the other links in
# this document (color 1155cc, single underline).
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("tin.nguyen@kyanon.digital", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found -and $rng.Hyperlinks.Count -eq 0) {
    $d.Hyperlinks.Add($rng, "mailto:tin.nguyen@kyanon.digital") | Out-Null

    # Re-find the (now-linked) text so we can apply the same direct formatting
    # used by the document's other hyperlinks.
    $rng2 = $d.Content.Duplicate
    $rng2.Find.Execute("tin.nguyen@kyanon.digital", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
    $rng2.Font.Color = 13391121   # RGB(0x11,0x55,0xCC) -> wdColor BGR long
    $rng2.Font.Underline = 1      # wdUnderlineSingle
}
